# Refresh crypto price (D) and 1h volume-change (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.130.34"
$ws.Range("E2").Value = "  +5.56%  "
$ws.Range("D3").Value = "2.341.31"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.63%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "2.701.45"
$ws.Range("E14").Value = "  +4.62%  "
$ws.Range("D15").Value = "2.340.65"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.48%  "
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "46.946.49"
$ws.Range("E18").Value = "  +5.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +16.68%  "
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0817"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "148.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("E36").Value = "  +5.12%  "
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0314"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.48%  "
$ws.Range("D45").Value = "1.831.49"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "75.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.195"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.17%  "
